$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: duplicate the formatting (and, for most columns, the exact
# values) of row 43 down into the brand-new row 44 BEFORE row 43 is
# touched, so that H44 naturally inherits the "IN PROGRESS" status
# that row 43 currently holds. (J is skipped on purpose - row 43 has
# no J cell yet at this point, and row 44 must not get one either.)
# ------------------------------------------------------------------
$ws.Range("A43:I43").Copy($ws.Range("A44:I44"))
$ws.Range("K43").Copy($ws.Range("K44"))

# ------------------------------------------------------------------
# Step 2: finalize row 43 - the Buy order that completed.
# ------------------------------------------------------------------
$ws.Range("H43").Value = "DONE"

# I43 gets the "finalized date" style/number-format used elsewhere in
# the same column (e.g. I41/I42) plus the actual finish time.
$ws.Range("I42").Copy($ws.Range("I43"))
$ws.Range("I43").Value = 42860.867835648147

# J43 is a brand-new cell (fee %), matching the plain/default style
# already used by the other J-column fee cells.
$ws.Range("J41").Copy($ws.Range("J43"))
$ws.Range("J43").Value = "0.56250000 XRP (0.15%)"

# ------------------------------------------------------------------
# Step 3: fill in row 44 - the new Sell order.
# ------------------------------------------------------------------
$ws.Range("A44").Value = 42861.347037037034

# B44 needs the same rich-text "Sell" run (leading spaces + red text)
# used elsewhere, so copy a real "Sell" cell wholesale.
$ws.Range("B42").Copy($ws.Range("B44"))

# C44/F44/G44 already match (same shared text as row 43) thanks to the
# Step 1 copy, so nothing further is required for those.

# D44 needs to stay a *text* value (it has leading spaces & embedded
# newlines) rather than being auto-parsed as a number. Writing it with
# a leading apostrophe forces text, then we re-apply the clean
# wrap-text/number format (without the quote-prefix flag) from D43.
$ws.Range("D44").Value = "'              0.09535348`n`n"
$ws.Range("D43").Copy()
$ws.Range("D44").PasteSpecial(-4122)

$ws.Range("E44").Value = "          0.107USDT"

# H44 already equals "IN PROGRESS" (copied in Step 1).

$excel.CutCopyMode = $false

# Fix the row height last, after all content/wrapping has been set,
# so Excel's auto-fit-on-entry doesn't override the explicit height.
$ws.Rows.Item(44).RowHeight = 14.25

# ------------------------------------------------------------------
# Step 4: sheet view - the user scrolled down and selected D45.
# ------------------------------------------------------------------
$ws.Range("D45").Select()
$excel.ActiveWindow.ScrollRow = 32
